$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 95-28 147th Place in Jamaica, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2026/01/permits-filed-for-95-28-147th-place-in-jamaica-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a 13-story residential building at 95-28 147th Place in <a href="https://newyorkyimby.com/neighborhoods/jamaica">Jamaica</a>, Queens. Located between 95th Avenue and 97th Avenue, the lot is near the Sutphin Boulevard–Archer Avenue–JFK Airport subway station, served by the E, J, and Z trains. Herman Jacob under East 181st Gardens LLC is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2026-01-29T11:30:10+00:00"
$ws.Range("E2").Value = "Thu, 29 Jan 2026 11:30:10 +0000"
$ws.Range("G2").Value = "YIMBY - Jamaica"
